$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "30EI7I"
$ws.Range("B7").Value = "Almohadilla + Chip Epson"
$ws.Range("C7").Value = "SC F500 F530 F531 F540 F551 F560 F570 F571 T2100 T2170 T3100 T3130 T3160 T3170 T3180 T3480 T5100 T5130 T5160 T5170 T5180"
$ws.Range("D7").Value = 35000
$ws.Range("E7").Value = 200000
$ws.Range("F7").Value = 7
$ws.Range("G7").Value = 1
$ws.Range("H7").Formula = "=(E7-D7)*G7"
$ws.Range("I7").Formula = "=D7*F7"
$ws.Range("J7").Value = 245000
